$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.224575996398926
$ws.Range("B1").Value = 2.778353214263916
$ws.Range("C1").Value = 2.99786639213562
$ws.Range("D1").Value = 3.592722415924072
$ws.Range("E1").Value = 1.78554356098175
